$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "2009年" row with full data); this shifts rows 3 and 4
# (2010年 / 2011年) up to become rows 2 and 3.
$ws.Rows.Item(2).Delete()
